# Final touches, go live
# Add a new "DESIGNATION" column (I) to the contacts template header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell - inherits the bold header style from the row, but set
# explicitly too for clarity/robustness.
$ws.Range("I1").Value = "DESIGNATION"
$ws.Range("I1").Font.Bold = $true

# Give the new column a sensible width (close to the author's final width).
$ws.Columns("I").ColumnWidth = 14.3

# Select column J (whole column), matching the cursor position the author
# ended up on after adding the new column.
$ws.Columns("J:J").Select() | Out-Null
